$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "TestSteps" (sheet1): insert a new accessibility-check step as row 3,
# pushing the existing steps (old rows 3-10) down to rows 4-11.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TestSteps")

# Shift existing data rows 3-10 down to rows 4-11 (carries values + styles).
$ws1.Range("A3:C10").Copy()
$ws1.Range("A4").PasteSpecial(-4104)

# The paste above only reached rows that existed in the prior used range, so
# the newly created row 11 didn't inherit formatting - fix it up explicitly.
$ws1.Range("A10:C10").Copy()
$ws1.Range("A11").PasteSpecial(-4122)

# New row 3: accessibility check step.
$ws1.Range("A3").Value = "checkAccessibility"
$ws1.Range("B3").Value = "AddNewProfile_GaapSourceObjProfile"

# C3:E3 stay empty but keep the same body style as the rest of the table.
$ws1.Range("C4").Copy()
$ws1.Range("C3:E3").PasteSpecial(-4122)
$ws1.Range("C3").ClearContents()

# F3 picks up a thin border (new cell style) to close off the row.
$ws1.Range("F3").Borders.LineStyle = 1

# Column widths widened to fit the new content.
$ws1.Columns.Item(1).ColumnWidth = 25.26953125
$ws1.Columns.Item(2).ColumnWidth = 33

# ---------------------------------------------------------------------------
# Sheet "TestData" (sheet2): bump the Category value used by the new test.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TestData")
$ws2.Range("J2").Value = " 16 - NACUBO REVENUES PROP "
$ws2.Columns.Item(10).ColumnWidth = 26.90625

# ---------------------------------------------------------------------------
# View state: TestData was the active/selected tab before, TestSteps is now.
# ---------------------------------------------------------------------------
$ws2.Range("E1").Select()
$ws1.Activate()
$ws1.Range("C5").Select()
